$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range('D2').Value = '26.320.91'
$ws.Range('E2').Value = '  -4.38%  '

$ws.Range('D3').Value = '1.754.58'
$ws.Range('E3').Value = '  -4.25%  '

$ws.Range('E4').Value = '  +0.10%  '

Set-TextValue 'D5' '1.002'
$ws.Range('E5').Value = '  +0.10%  '

Set-TextValue 'D6' '303.07'
$ws.Range('E6').Value = '  -2.88%  '

Set-TextValue 'D7' '0.4285'
$ws.Range('E7').Value = '  +0.28%  '

$ws.Range('E8').Value = '  -1.67%  '

Set-TextValue 'D9' '0.07014'
$ws.Range('E9').Value = '  -3.41%  '

Set-TextValue 'D10' '0.8269'
$ws.Range('E10').Value = '  -4.13%  '

Set-TextValue 'D11' '20.03'
$ws.Range('E11').Value = '  -2.78%  '

$ws.Range('D12').Value = '1.738.25'
$ws.Range('E12').Value = '  -4.58%  '

Set-TextValue 'D13' '5.184'
$ws.Range('E13').Value = '  -3.91%  '

Set-TextValue 'D14' '6.303'
$ws.Range('E14').Value = '  -3.17%  '

Set-TextValue 'D15' '0.06790'
$ws.Range('E15').Value = '  -2.10%  '

Set-TextValue 'D16' '1.007'
$ws.Range('E16').Value = '  +0.56%  '

Set-TextValue 'D17' '78.67'
$ws.Range('E17').Value = '  -2.39%  '

Set-TextValue 'D18' '0.000008621'
$ws.Range('E18').Value = '  -3.15%  '

Set-TextValue 'D19' '1.004'
$ws.Range('E19').Value = '  +0.34%  '

Set-TextValue 'D20' '14.81'
$ws.Range('E20').Value = '  -3.75%  '

$ws.Range('D21').Value = '26.351.03'
$ws.Range('E21').Value = '  -4.04%  '

Set-TextValue 'D22' '4.952'
$ws.Range('E22').Value = '  -3.76%  '

Set-TextValue 'D23' '11.06'
$ws.Range('E23').Value = '  +2.31%  '

$ws.Range('D24').Value = '1.978.36'
$ws.Range('E24').Value = '  -3.82%  '

Set-TextValue 'D25' '1.902'
$ws.Range('E25').Value = '  -4.39%  '

Set-TextValue 'D26' '151.70'
$ws.Range('E26').Value = '  -1.83%  '

Set-TextValue 'D27' '18.03'
$ws.Range('E27').Value = '  -4.23%  '

Set-TextValue 'D28' '114.65'
$ws.Range('E28').Value = '  +0.30%  '

Set-TextValue 'D29' '4.998'
$ws.Range('E29').Value = '  -1.97%  '

Set-TextValue 'D30' '1.637'
$ws.Range('E30').Value = '  -9.69%  '

Set-TextValue 'D31' '0.08902'
$ws.Range('E31').Value = '  +0.57%  '

Set-TextValue 'D32' '0.7124'
$ws.Range('E32').Value = '  -4.14%  '

Set-TextValue 'D33' '4.274'
$ws.Range('E33').Value = '  -5.63%  '

Set-TextValue 'D34' '1.085'
$ws.Range('E34').Value = '  -3.86%  '

Set-TextValue 'D35' '1.002'
$ws.Range('E35').Value = '  +0.13%  '

Set-TextValue 'D36' '2.740'
$ws.Range('E36').Value = '  -8.24%  '

$ws.Range('E37').Value = '  -2.74%  '

Set-TextValue 'D38' '0.05059'
$ws.Range('E38').Value = '  -4.58%  '

Set-TextValue 'D39' '0.01871'
$ws.Range('E39').Value = '  -3.11%  '

$ws.Range('B40').Value = 'Algorand'
$ws.Range('C40').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 'D40' '0.1591'
$ws.Range('E40').Value = '  -3.92%  '

$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue 'D41' '0.4848'
$ws.Range('E41').Value = '  -4.32%  '

$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D42' '6.112'
$ws.Range('E42').Value = '  -5.46%  '

$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D43' '2.461'
$ws.Range('E43').Value = '  -12.00%  '

$ws.Range('B44').Value = 'Aptos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D44' '7.848'
$ws.Range('E44').Value = '  -5.41%  '

$ws.Range('B45').Value = 'Quant'
$ws.Range('C45').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue 'D45' '104.19'
$ws.Range('E45').Value = '  -1.00%  '

$ws.Range('B46').Value = 'PaxDollar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue 'D46' '1.002'
$ws.Range('E46').Value = '  +0.21%  '

Set-TextValue 'D47' '9.949'
$ws.Range('E47').Value = '  -3.98%  '

Set-TextValue 'D48' '0.06172'
$ws.Range('E48').Value = '  -4.73%  '

Set-TextValue 'D49' '0.4437'
$ws.Range('E49').Value = '  -4.95%  '

Set-TextValue 'D50' '1.557'
$ws.Range('E50').Value = '  -3.26%  '

Set-TextValue 'D51' '1.697'
$ws.Range('E51').Value = '  -0.93%  '
